$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value (45182 -> 2023-09-13) for
# every data row (2..115). The commit updates that "changed on" date to
# 45184 (2023-09-15) for all of them.
for ($r = 2; $r -le 115; $r++) {
    $ws.Cells.Item($r, 3).Value = 45184
}
